$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndex = @{ B = 2; C = 3; D = 4; E = 5 }

$changes = @{
    2 = @{ D="93.675.52"; E="  +1.28%  " }
    3 = @{ D="3.490.14"; E="  +4.93%  " }
    4 = @{ E="  +0.04%  " }
    5 = @{ D="'235.34"; E="  +3.15%  " }
    6 = @{ D="'625.26"; E="  +0.60%  " }
    7 = @{ D="'1.43"; E="  +7.41%  " }
    8 = @{ D="'0.392"; E="  +4.84%  " }
    9 = @{ E="  +0.03%  " }
    10 = @{ D="'0.996"; E="  +10.14%  " }
    11 = @{ D="3.481.32"; E="  +4.70%  " }
    12 = @{ D="'42.80"; E="  +2.51%  " }
    13 = @{ D="'0.200"; E="  +5.67%  " }
    14 = @{ D="'6.26"; E="  +5.81%  " }
    15 = @{ D="4.151.06"; E="  +5.09%  " }
    16 = @{ D="93.456.76"; E="  +1.55%  " }
    17 = @{ D="'0.0000249"; E="  +4.36%  " }
    18 = @{ D="'8.33"; E="  +6.24%  " }
    19 = @{ D="3.496.16"; E="  +5.41%  " }
    20 = @{ B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="'12.43"; E="  +14.91%  " }
    21 = @{ B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="'17.97"; E="  +7.38%  " }
    22 = @{ D="'0.497"; E="  +14.61%  " }
    23 = @{ B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="'517.86"; E="  +6.76%  " }
    24 = @{ B="SuiNetwork"; C="https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"; D="'3.38"; E="  +4.08%  " }
    25 = @{ D="'6.72"; E="  +9.96%  " }
    26 = @{ D="'0.0000182"; E="  +2.00%  " }
    27 = @{ D="'95.48"; E="  +7.32%  " }
    28 = @{ D="'12.23"; E="  +6.97%  " }
    29 = @{ D="3.672.10"; E="  +5.05%  " }
    30 = @{ D="'2.95"; E="  +14.09%  " }
    31 = @{ B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="'11.36"; E="  +3.37%  " }
    32 = @{ B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="'1.00"; E="  +0.04%  " }
    33 = @{ D="'0.137"; E="  +3.98%  " }
    34 = @{ D="'1.01"; E="  +0.49%  " }
    35 = @{ D="'0.178"; E="  +5.94%  " }
    36 = @{ D="'29.81"; E="  +6.67%  " }
    37 = @{ D="'0.557"; E="  +7.91%  " }
    38 = @{ D="'567.82"; E="  +9.49%  " }
    39 = @{ E="  +7.74%  " }
    40 = @{ D="'7.53"; E="  +4.24%  " }
    41 = @{ E="  -0.02%  " }
    42 = @{ D="'0.922"; E="  +6.15%  " }
    43 = @{ D="'0.148"; E="  +2.47%  " }
    44 = @{ D="'23.75"; E="  -0.95%  " }
    45 = @{ B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="'1.69"; E="  +3.21%  " }
    46 = @{ B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="'0.0417"; E="  +7.77%  " }
    47 = @{ D="'3.55"; E="  -0.66%  " }
    48 = @{ D="'5.45"; E="  +3.34%  " }
    49 = @{ B="Stacks"; C="https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D="'2.15"; E="  +3.20%  " }
    50 = @{ B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="'53.44"; E="  +3.11%  " }
    51 = @{ D="'8.13"; E="  +3.97%  " }
}

foreach ($row in $changes.Keys) {
    $rowData = $changes[$row]
    foreach ($col in $rowData.Keys) {
        $colNum = $colIndex[$col]
        $ws.Cells.Item([int]$row, $colNum).Value = $rowData[$col]
    }
}

Write-Output "Applied $($changes.Count) row updates"
